$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Jacksonville altimeter weather question: swap Correct/Incorrect answers
$ws.Range("D5").Value = "Lower"
$ws.Range("E5").Value = "Higher"

# Copy formatting from template rows down to the new rows
$ws.Range("A2:H2").Copy()
$ws.Range("A106:H118").PasteSpecial(-4122)
$ws.Range("A2:B2").Copy()
$ws.Range("A119:B119").PasteSpecial(-4122)

# Add new engines questions
$ws.Range("A106").Value = "engines"
$ws.Range("B106").Value = 5
$ws.Range("C106").Value = "The power control system refers to:"
$ws.Range("D106").Value = "The hydraulic system that provides pressure only for the flight controls"
$ws.Range("E106").Value = "The part of the AC power system that powers essential equipment for flight"
$ws.Range("F106").Value = "The system which controls engine power output, specifically the Power Control Lever, FCUs, and pressurizing and dump valves"
$ws.Range("G106").Value = "The part of the lubrication which reacts to engine power output changes"

$ws.Range("A107").Value = "engines"
$ws.Range("B107").Value = 5
$ws.Range("C107").Value = "The Accumulator works in tandem with the _______ to maintain system pressure during shutdown"
$ws.Range("D107").Value = "Check Valve"
$ws.Range("E107").Value = "Pressure Relief Valve"
$ws.Range("F107").Value = "Pressure Regulator Valve"
$ws.Range("G107").Value = "Pressurizing Valve"

$ws.Range("A108").Value = "engines"
$ws.Range("B108").Value = 1
$ws.Range("C108").Value = "Who discovered the inverse relationship between static and dynamic pressure in a closed system?"
$ws.Range("D108").Value = "Bernoulli"
$ws.Range("E108").Value = "Pascal"
$ws.Range("F108").Value = "Otto"
$ws.Range("G108").Value = "Brayton"

$ws.Range("A109").Value = "engines"
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = "As Engine RPM increases, thrust ________. This effect gets _______ pronounced at higher RPM"
$ws.Range("D109").Value = "Increases, more"
$ws.Range("E109").Value = "Decreases, less"
$ws.Range("F109").Value = "Decreases, more"
$ws.Range("G109").Value = "Increases, less"

$ws.Range("A110").Value = "engines"
$ws.Range("B110").Value = 1
$ws.Range("C110").Value = "As altitude increases, thrust ______"
$ws.Range("D110").Value = "Decreases slowly then faster as temperature initially decreases then levels off, while pressure constantly decreases"
$ws.Range("E110").Value = "Decreases at a constant rate as temperature decreases and pressure decreases at a faster rate"
$ws.Range("F110").Value = "Increases as temperature decreases at a faster rate than pressure decreases"
$ws.Range("G110").Value = "Increases initially then decreases as temperature initially decreases then levels off, while pressure constantly decreases"

$ws.Range("A111").Value = "engines"
$ws.Range("B111").Value = 1
$ws.Range("C111").Value = "Because of ram effect, as airspeed increases, density ______ and thrust ______"
$ws.Range("D111").Value = "Increases, increases"
$ws.Range("E111").Value = "Increases, decreases"
$ws.Range("F111").Value = "Decreases, increases"
$ws.Range("G111").Value = "Decreases, decreases"

$ws.Range("A112").Value = "engines"
$ws.Range("B112").Value = 1
$ws.Range("C112").Value = "As supersonic airflow goes throw a convergent passage, total pressure"
$ws.Range("D112").Value = "Remains constant"
$ws.Range("E112").Value = "Increases"
$ws.Range("F112").Value = "Decreases"
$ws.Range("G112").Value = "Decreases then Increases"

$ws.Range("A113").Value = "engines "
$ws.Range("B113").Value = 1
$ws.Range("C113").Value = "Airflow through a variable geometry inlet duct first go through a _______ passage then a _______ passage in order to maximize _________"
$ws.Range("D113").Value = "Convergent, divergent, pressure"
$ws.Range("E113").Value = "Convergent, divergent, velocity"
$ws.Range("F113").Value = "Divergent, convergent, pressure"
$ws.Range("G113").Value = "Divergent, convergent, velocity"

$ws.Range("A114").Value = "engines "
$ws.Range("B114").Value = 4
$ws.Range("C114").Value = "In a reciprocating engine "
$ws.Range("D114").Value = "Directly to the crankshaft"
$ws.Range("E114").Value = "Directly to the camshaft"
$ws.Range("F114").Value = "Through the constant speed drive"
$ws.Range("G114").Value = "Through the governor"

$ws.Range("A115").Value = "engines"
$ws.Range("B115").Value = 4
$ws.Range("C115").Value = "In a dual spool axial compressor GTE, what drives the propeller?"
$ws.Range("D115").Value = "The low pressure compressor"
$ws.Range("E115").Value = "The high pressure turbine"
$ws.Range("F115").Value = "The high pressure compressor"
$ws.Range("G115").Value = "Exhaust gases"

$ws.Range("A116").Value = "engines"
$ws.Range("B116").Value = 4
$ws.Range("C116").Value = "A gas turbine engine is powered by the _____ cycle which occurs _____"
$ws.Range("D116").Value = "Brayton, simultaneously"
$ws.Range("E116").Value = "Brayton, sequentially"
$ws.Range("F116").Value = "Otto, simultaneously"
$ws.Range("G116").Value = "Otto, sequentially"

$ws.Range("A117").Value = "engines"
$ws.Range("B117").Value = 4
$ws.Range("C117").Value = "A reciprocating engine is powered by the _____ cycle which occurs _____"
$ws.Range("D117").Value = "Otto, sequentially"
$ws.Range("E117").Value = "Otto, simultaneously"
$ws.Range("F117").Value = "Brayton, sequentially"
$ws.Range("G117").Value = "Brayton, simultaneously"

$ws.Range("A118").Value = "engines"
$ws.Range("B118").Value = 1
$ws.Range("C118").Value = "Torquemeters measure ______ and are found in _______"
$ws.Range("D118").Value = "Shaft horsepower, turbo props and turbo fans"
$ws.Range("E118").Value = "Torque, turbo props and turbo fans"
$ws.Range("F118").Value = "Foot pounds, turbo props and turbo fans"
$ws.Range("G118").Value = "Torque, turbo fans"

$ws.Range("A119").Value = "engines"
$ws.Range("B119").Value = 1

# Extend the ID (row number) formula to the new rows
$ws.Range("H106:H118").Formula = "=ROW()"